$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 416793
$ws.Range("C2").Value = "Sanskriti The Gurukul"
$ws.Range("D2").Value = "Guwahati"

$ws.Range("E2").Select()
